$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block 1: Actions_3.p (rows 28-38)
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Actions_3.p"
$ws.Range("B28").Value = "Feature Integration"
$ws.Range("C28").Value = "No"
$ws.Range("D28").Value = "Yes"
$ws.Range("E28").Value = "No"
$ws.Range("F28").Value = "1.7.2. Transition Declaration: push transition"
$ws.Range("G28").Value = "Based on Actions_1.p, but with push and pop added"

$ws.Range("F29").Value = "1.2.2. Real and model machines"
$ws.Range("F30").Value = "1.8.2. Do declaration: action on named function"
$ws.Range("F31").Value = "2.2. Dynamic creation of machines using ""new"""
$ws.Range("F32").Value = "2.3. ""raise"" stmt"
$ws.Range("F33").Value = "2.4. ""send"" stmt"
$ws.Range("F34").Value = "2.8. Assertions"
$ws.Range("F35").Value = "2.9.1. Assign on bool"
$ws.Range("F36").Value = "3.3.1. ""payload"" primitive expression"
$ws.Range("F37").Value = "3.3.6. Cast operator ""as"""
$ws.Range("F38").Value = "4.3. Passing variables as payloads: variable of ""ghost machine"" type"

# ---------------------------------------------------------------------------
# Block 2: Actions_4.p (rows 40-50) -- row 39 intentionally left blank
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "Actions_4.p"
$ws.Range("B40").Value = "Feature Integration"
$ws.Range("C40").Value = "No"
$ws.Range("D40").Value = "Yes"
$ws.Range("E40").Value = "No"
$ws.Range("F40").Value = "1.7.2. Transition Declaration: push transition"
$ws.Range("G40").Value = "Similar to Actions_3.p, but with two asserts in a row: the 1st assert passes, the 2nd assert fails"

$ws.Range("F41").Value = "1.2.2. Real and model machines"
$ws.Range("F42").Value = "1.8.2. Do declaration: action on named function"
$ws.Range("F43").Value = "2.2. Dynamic creation of machines using ""new"""
$ws.Range("F44").Value = "2.3. ""raise"" stmt"
$ws.Range("F45").Value = "2.4. ""send"" stmt"
$ws.Range("F46").Value = "2.8. Assertions: two asserts in a row"
$ws.Range("F47").Value = "2.9.1. Assign on bool"
$ws.Range("F48").Value = "3.3.1. ""payload"" primitive expression"
$ws.Range("F49").Value = "3.3.6. Cast operator ""as"""
$ws.Range("F50").Value = "4.3. Passing variables as payloads: variable of ""ghost machine"" type"

# Bold the "two asserts in a row" portion of F46 (rich text run)
$chars = $ws.Range("F46").Characters(18, 21)
$chars.Font.Bold = $true

# ---------------------------------------------------------------------------
# Block 3: Actions_5.p / BangaloreToRedmond / Call_Exit_* / CallStatement_1
# ---------------------------------------------------------------------------
$ws.Range("A51").Value = "Actions_5.p"
$ws.Range("B51").Value = "Feature Integration"
$ws.Range("F51").Value = "TODO"

$ws.Range("A52").Value = "BangaloreToRedmond"
$ws.Range("B52").Value = "Feature Integration"
$ws.Range("C52").Value = "No"
$ws.Range("D52").Value = "Yes"
$ws.Range("E52").Value = "No"
$ws.Range("F52").Value = "TODO"

$ws.Range("A53").Value = "Call_Exit_1"
$ws.Range("B53").Value = "Feature Integration"
$ws.Range("C53").Value = "No"
$ws.Range("D53").Value = "Yes"
$ws.Range("E53").Value = "No"
$ws.Range("F53").Value = "TODO"
$ws.Range("G53").Value = "Identical to BangaloreToRedmond, but using ""with"": on default goto TakeBus with { RemoteCheckIn = true; };"

$ws.Range("A54").Value = "Call_Exit_2"
$ws.Range("B54").Value = "Feature Integration"
$ws.Range("C54").Value = "No"
$ws.Range("D54").Value = "Yes"
$ws.Range("E54").Value = "No"
$ws.Range("F54").Value = "TODO"

$ws.Range("A55").Value = "Call_Exit_3"
$ws.Range("B55").Value = "Feature Integration"
$ws.Range("C55").Value = "No"
$ws.Range("D55").Value = "Yes"
$ws.Range("E55").Value = "No"
$ws.Range("F55").Value = "TODO"

$ws.Range("A56").Value = "CallStatement_1"
$ws.Range("B56").Value = "Feature Integration"
$ws.Range("C56").Value = "No"
$ws.Range("D56").Value = "Yes"
$ws.Range("E56").Value = "No"
$ws.Range("F56").Value = "TODO"

# ---------------------------------------------------------------------------
# Formatting: yellow "RegressionTests" rows, matching rows 6/18 etc.
# ---------------------------------------------------------------------------
$ws.Range("A28:E28").Interior.Color = 65535
$ws.Range("G28").Interior.Color = 65535
$ws.Range("F28").Interior.Color = 65535
$ws.Range("F28").Font.Bold = $true

$ws.Range("F29").Interior.Color = 65535
$ws.Range("F29").Interior.ColorIndex = -4142

$ws.Range("A40:E40").Interior.Color = 65535
$ws.Range("G40").Interior.Color = 65535
$ws.Range("F40").Interior.Color = 65535
$ws.Range("F40").Font.Name = "Calibri"

$ws.Range("F41").Interior.Color = 65535
$ws.Range("F41").Interior.ColorIndex = -4142

$ws.Range("A51:B51").Interior.Color = 65535
$ws.Range("F51").Interior.Color = 65535

$ws.Range("A52:F52").Interior.Color = 65535
$ws.Range("A53:G53").Interior.Color = 65535
$ws.Range("A54:F54").Interior.Color = 65535
$ws.Range("A55:F55").Interior.Color = 65535
$ws.Range("A56:F56").Interior.Color = 65535

# ---------------------------------------------------------------------------
# View state: selection on the last newly-populated cell
# ---------------------------------------------------------------------------
$ws.Range("G53").Select()
